$d = $word.ActiveDocument

# 1. Insert a new agenda topic after "Introduce the meeting" (i.e. right
#    before "Feedback on Hand-In document"), explaining why there is not
#    much demo-able production.
$introPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Introduce the meeting*") {
        $introPara = $p
        break
    }
}
$introPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($introPara.Index + 1)
$newPara.Range.Text = "Explanation on why there is not much demo-able production – 10 mins"

# 2. "Feedback on Hand-In document" now only needs 5 minutes instead of 10.
$d.Content.Find.Execute("Feedback on Hand-In document – 10 mins", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Feedback on Hand-In document – 5 mins", 2) | Out-Null

# 3. Drop the parenthetical detail from the sprint deliverables topic and
#    shorten its slot to 10 minutes.
$d.Content.Find.Execute("Talk about sprint deliverables (Payment Service Continued, Many-to-Many Broadcast) – 20 mins", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Talk about sprint deliverables  – 10 mins", 2) | Out-Null

# 4. "Questions from client" shrinks from 5 to 2 minutes.
$d.Content.Find.Execute("Questions from client – 5 mins", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Questions from client – 2 mins", 2) | Out-Null

# 5. "Questions from developers (company)" stays the same (no change needed).

# 6. "Next steps discussion with teacher" grows from 5 to 10 minutes.
$d.Content.Find.Execute("Next steps discussion with teacher – 5 mins", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Next steps discussion with teacher – 10 mins", 2) | Out-Null

# 7. Remove the now-redundant standalone "Questions from teacher" topic.
$removePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Questions from teacher*") {
        $removePara = $p
        break
    }
}
$removePara.Range.Delete()

# 8. "Questions from students" grows from 2 to 10 minutes.
$d.Content.Find.Execute("Questions from students – 2 mins", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Questions from students – 10 mins", 2) | Out-Null
